# Mark task "refactor all panels to use proper events instead of delegate
# properties" (Id 100) as Done: move its row from the Active sheet to the
# top of the Inactive sheet, updating its Status to Done and stamping a
# Done date equal to its Created date.

$wb = $excel.ActiveWorkbook
$active = $wb.Worksheets.Item("Active")
$inactive = $wb.Worksheets.Item("Inactive")

# Capture the row-2 values from the Active sheet before removing it.
$id = $active.Range("A2").Value()
$title = $active.Range("B2").Value()
$category = $active.Range("D2").Value()
$created = $active.Range("E2").Value()

# Remove the task from the Active (Todo) sheet; remaining rows shift up.
$active.Rows.Item(2).Delete()

# Insert a new row at the top of the Inactive (Done) sheet for this task.
$inactive.Rows.Item(2).Insert()

# Insert() copies the formatting of the row above (the bold header row) onto
# the new row; re-normalize it to match the plain data rows below instead.
$inactive.Range("A3:F3").Copy()
$inactive.Range("A2:F2").PasteSpecial(-4122)

$inactive.Range("A2").Value = $id
$inactive.Range("B2").Value = $title
$inactive.Range("C2").Value = "Done"
$inactive.Range("D2").Value = $category
# Dates are stored as plain text in this sheet, not real date serials; a
# leading apostrophe keeps the assignment as text instead of letting Excel
# auto-convert the "m/d/yyyy"-looking string into a date value. Re-apply the
# plain data-row style afterwards so the forced-text quote prefix doesn't
# stick around as a formatting difference.
$inactive.Range("E2").Value = "'" + $created
$inactive.Range("F2").Value = "'" + $created
$inactive.Range("E2").Style = $inactive.Range("E3").Style
$inactive.Range("F2").Style = $inactive.Range("F3").Style
